# Apply the dataset update: add new feature columns to the Host and Dopant
# tables, then leave the workbook with Dopant as the active sheet (matching
# the author's final view state).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Host sheet (Table1): add Host_d_band_filling, Host_WSR, Host_surface_energy
# ---------------------------------------------------------------------
$wsHost = $wb.Worksheets.Item("Host")
$tblHost = $wsHost.ListObjects.Item("Table1")

$colDBand = $tblHost.ListColumns.Add()
$colDBand.Range.Cells(1,1).Value = "Host_d_band_filling"

$colWSR = $tblHost.ListColumns.Add()
$colWSR.Range.Cells(1,1).Value = "Host_WSR"

$colSurfE = $tblHost.ListColumns.Add()
$colSurfE.Range.Cells(1,1).Value = "Host_surface_energy"

# Data rows (Cu, Ag, Au) -> Host_At_No 29, 47, 79
$hostDBand = @(1, 1, 1)
$hostWSR = @(2.67, 3.01, 3)
$hostSurfE = @(1.8075000000000001, 1.248, 1.5029999999999999)

for ($i = 0; $i -lt 3; $i++) {
    $r = $i + 2
    $wsHost.Range("K$r").Value = $hostDBand[$i]
    $wsHost.Range("K$r").NumberFormat = "0.00"
    $wsHost.Range("L$r").Value = $hostWSR[$i]
    $wsHost.Range("M$r").Value = $hostSurfE[$i]
}

# ---------------------------------------------------------------------
# Dopant sheet (Table3): add Dopant_d_band_filling, Dopant_surface_energy,
# Dopant_WSR, Dopant_d_band_centre
# ---------------------------------------------------------------------
$wsDopant = $wb.Worksheets.Item("Dopant")
$tblDopant = $wsDopant.ListObjects.Item("Table3")

$colDopDBand = $tblDopant.ListColumns.Add()
$colDopDBand.Range.Cells(1,1).Value = "Dopant_d_band_filling"

$colDopSurfE = $tblDopant.ListColumns.Add()
$colDopSurfE.Range.Cells(1,1).Value = "Dopant_surface_energy"

$colDopWSR = $tblDopant.ListColumns.Add()
$colDopWSR.Range.Cells(1,1).Value = "Dopant_WSR"

$colDopDCentre = $tblDopant.ListColumns.Add()
$colDopDCentre.Range.Cells(1,1).Value = "Dopant_d_band_centre"

# Rows 2-21, columns J (d_band_filling), K (surface_energy), L (WSR), M (d_band_centre)
$dopJ = @(0.2, 0.3, 0.4, 0.7, 0.8, 0.9, 1, 0.2, 0.3, 0.4, 0.5, 0.7, 0.8, 0.9, 1, 0.5, 0.6, 0.8, 0.9, 1)
$dopK = @(1.2749999999999999, 2.0455000000000001, 2.5859999999999999, 2.4460000000000002, 2.536, 2.415, 1.8075000000000001, 1.125, 1.9544999999999999, 2.6675, 2.9535, 3.0465, 2.6795, 2.0265, 1.248, 3.47, 3.613, 3.024, 2.4820000000000002, 1.5029999999999999)
$dopL = @(3.43, 3.05, 2.82, 2.66, 2.62, 2.6, 2.67, 3.76, 3.35, 3.07, 2.99, 2.79, 2.81, 2.87, 3.01, 2.95, 2.87, 2.84, 2.9, 3)
$dopM = @(2, 1.5, 1.06, -0.92, -1.17, -1.29, -2.67, 3, 1.95, 1.41, -0.6, -1.41, -1.73, -1.83, -4.3, 0.77, -0.51, -2.11, -2.25, -3.56)

for ($i = 0; $i -lt 20; $i++) {
    $r = $i + 2
    $wsDopant.Range("J$r").Value = $dopJ[$i]
    $wsDopant.Range("K$r").Value = $dopK[$i]
    $wsDopant.Range("L$r").Value = $dopL[$i]
    $wsDopant.Range("M$r").Value = $dopM[$i]
}

# ---------------------------------------------------------------------
# View / selection state: Dopant becomes the active sheet; Host keeps a
# selection on its new data; Final loses its "active" flag.
# ---------------------------------------------------------------------
$wsHost.Activate()
$wsHost.Range("M3").Select()

$wsDopant.Activate()
$wsDopant.Range("M10").Select()
